$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.323867321014404
$ws.Range("B1").Value = 2.632214546203613
$ws.Range("C1").Value = 2.587695598602295
$ws.Range("D1").Value = 1.949337482452393
$ws.Range("E1").Value = 0.5532294511795044
